$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 20: A20 (date), B20 (task text) ---
$ws.Range("A20").Value = 45164
$ws.Range("B20").Value = "cluster PCA results of images"

# Copy formatting from row 19 (A19:B19) onto the new row 20 cells so that
# style indices / borders match the existing table look (s=13 / s=14).
$null = $ws.Range("A19:B19").Copy()
$null = $ws.Range("A20:B20").PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = 18

# --- New row 27: D27 (question text) ---
$ws.Range("D27").Value = "why are cluster sizes of PCA results imbalanced?"

# Copy formatting from D26 onto D27 (s=8).
$null = $ws.Range("D26").Copy()
$null = $ws.Range("D27").PasteSpecial(-4122)
$ws.Rows.Item(27).RowHeight = 17

# --- Update selection / active cell / scroll position to reflect the
#     newly added row (sheetView topLeftCell="A19", selection B20) ---
$null = $ws.Range("A20:B20").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

$excel.CutCopyMode = $false
